# Weekly update: insert the newest "Betarraga" price-report week (Primera /
# Segunda quality rows) ahead of the existing history, shifting the older
# weeks down by two rows (the sheet is ordered most-recent-week-first after
# the first row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 150:169 down to 152:171, duplicating formatting
# (incl. the date-number style on column D) the same way Excel's own
# "Insert Copied/Blank Cells" does.
$ws.Rows("150:151").Insert()

# Populate the two freshly-inserted rows with this week's figures.
# Row 150 - Calidad "Primera"
$ws.Range("A150").Value = 8
$ws.Range("B150").Value = "Terminal La Palmera de La Serena"
$ws.Range("C150").Value = "Coquimbo"
$ws.Range("D150").Value = 44474
$ws.Range("E150").Value = 4
$ws.Range("F150").Value = 100114014
$ws.Range("G150").Value = "Betarraga"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 3000
$ws.Range("K150").Value = 450
$ws.Range("L150").Value = 500
$ws.Range("M150").Value = 475
$ws.Range("N150").Value = "`$/paquete 3 unidades"
$ws.Range("O150").Value = "Provincia del Elquí"
$ws.Range("P150").Value = 158
$ws.Range("Q150").Value = 3
$ws.Range("R150").Value = "Hortaliza"

# Row 151 - Calidad "Segunda"
$ws.Range("A151").Value = 8
$ws.Range("B151").Value = "Terminal La Palmera de La Serena"
$ws.Range("C151").Value = "Coquimbo"
$ws.Range("D151").Value = 44474
$ws.Range("E151").Value = 4
$ws.Range("F151").Value = 100114014
$ws.Range("G151").Value = "Betarraga"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Segunda"
$ws.Range("J151").Value = 1360
$ws.Range("K151").Value = 350
$ws.Range("L151").Value = 400
$ws.Range("M151").Value = 375
$ws.Range("N151").Value = "`$/paquete 3 unidades"
$ws.Range("O151").Value = "Provincia del Elquí"
$ws.Range("P151").Value = 125
$ws.Range("Q151").Value = 3
$ws.Range("R151").Value = "Hortaliza"
